# Update "Sprint Plan & Retorspective.xlsx" - Retrospective sheet:
#  - restyle the Sprint #1 / Sprint #2 label cells (top-aligned banding)
#  - fill in the Sprint #2 retrospective answers (row 5)
#  - add a new Sprint #3 retrospective row (row 6)
#  - fix up the sheet selection / active tab to match the saved state
#
# NOTE: every new cell style below is first assembled on a scratch cell
# (Z100) using Copy/PasteSpecial(formats) so that each distinct alignment
# combination is only ever committed to the stylesheet ONCE (building it up
# property-by-property directly on the destination cell would otherwise
# record each intermediate combination as its own separate cellXf).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$scratch = $ws2.Range("Z100")

# --- Style used for B5:D5 (Sprint #2 answer cells): ------------------------
# fillId=5 banding (same fill as the Plan sheet's Sprint #2 blocks),
# wrapped text, top vertical alignment, no horizontal alignment.
$ws1.Range("A12").Copy()
$scratch.PasteSpecial(-4122)
$scratch.WrapText = $true
$scratch.VerticalAlignment = -4160
$scratch.Copy()
$ws2.Range("B5:D5").PasteSpecial(-4122)
$scratch.Clear()

$ws2.Range("B5").Value = "Members who were able to work on their tasks managed to get the majority if not all their work doe within in the sprint"
$ws2.Range("C5").Value = "All members can do their best to communicate when unable to do work so others may make up the difference before the sprint finishes."
$ws2.Range("D5").Value = "Task estimates will be done with a points system to gauge the best estimates for each task. The team's velocity will also be determined next sprint."

# --- Style used for A4 & A6 ("Sprint #1" / "Sprint #3" labels): ------------
# Keeps the existing fillId=2 banding, adds top vertical alignment.
$ws2.Range("A4").Copy()
$scratch.PasteSpecial(-4122)
$scratch.VerticalAlignment = -4160
$scratch.Copy()
$ws2.Range("A4").PasteSpecial(-4122)
$scratch.Clear()

# --- Style used for A5 ("Sprint #2" label): ---------------------------------
# fillId=5 banding (matching B5:D5 above) + top vertical alignment.
$ws1.Range("A12").Copy()
$scratch.PasteSpecial(-4122)
$scratch.VerticalAlignment = -4160
$scratch.Copy()
$ws2.Range("A5").PasteSpecial(-4122)
$scratch.Clear()

# --- New row 6: Sprint #3 retrospective answers -----------------------------
# A6 reuses the same style as A4 (fillId=2 + top alignment).
$ws2.Range("A4").Copy()
$ws2.Range("A6").PasteSpecial(-4122)
$ws2.Range("A6").Value = "Sprint #3"

# B6:D6 reuse the fillId=2 / wrap / top banding already used on row 4.
$ws2.Range("B4:D4").Copy()
$ws2.Range("B6:D6").PasteSpecial(-4122)

$ws2.Range("C6").Value = "members could try to find some more motivation towards the project as there is only 3 more sprints to go"
$ws2.Range("B6").Value = "Not much really, team motivation was down due to the repetitive nature of the project"
$ws2.Range("D6").Value = "the group will try to finish up older tasks in the next sprint, delaying the final event to sprint 5"

# --- Row heights -------------------------------------------------------------
$ws2.Rows.Item(5).RowHeight = 75
$ws2.Rows.Item(6).RowHeight = 60

# --- Sheet/view selection ----------------------------------------------------
$ws2.Activate()
$ws2.Range("H5").Select()
